$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.636417666666667
$ws.Range("H2").Value = 22.909253
$ws.Range("I2").Value = 0.108532481296676
$ws.Range("J2").Value = 0.108532481296676
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.345835
$ws.Range("N2").Value = 91.03750500000001
$ws.Range("O2").Value = 0.8527782452855476
$ws.Range("P2").Value = 0.8527782452855475
$ws.Range("Q2").Value = 231.7334705037517
$ws.Range("R2").Value = 2085.601234533765
$ws.Range("S2").Value = 0.09255413895666585
$ws.Range("T2").Value = 0.09255413895666585

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.636417666666667
$ws.Range("H3").Value = 22.909253
$ws.Range("I3").Value = 0.108532481296676
$ws.Range("J3").Value = 0.108532481296676
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.224107666666666
$ws.Range("N3").Value = 9.672322999999999
$ws.Range("O3").Value = 0.09060383010029813
$ws.Range("P3").Value = 0.09060383010029811
$ws.Range("Q3").Value = 24.62063274496877
$ws.Range("R3").Value = 221.585694704719
$ws.Range("S3").Value = 0.009833458495767814
$ws.Range("T3").Value = 0.009833458495767814

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.636417666666667
$ws.Range("H4").Value = 22.909253
$ws.Range("I4").Value = 0.108532481296676
$ws.Range("J4").Value = 0.108532481296676
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.014730333333333
$ws.Range("N4").Value = 6.044191
$ws.Range("O4").Value = 0.05661792461415433
$ws.Range("P4").Value = 0.05661792461415433
$ws.Range("Q4").Value = 15.38532231103589
$ws.Range("R4").Value = 138.467900799323
$ws.Range("S4").Value = 0.006144883844242315
$ws.Range("T4").Value = 0.006144883844242316

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.103385
$ws.Range("H5").Value = 45.31015499999999
$ws.Range("I5").Value = 0.214656652056136
$ws.Range("J5").Value = 0.214656652056136
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.345835
$ws.Range("N5").Value = 91.03750500000001
$ws.Range("O5").Value = 0.8527782452855476
$ws.Range("P5").Value = 0.8527782452855475
$ws.Range("Q5").Value = 458.324829151475
$ws.Range("R5").Value = 4124.923462363275
$ws.Range("S5").Value = 0.183054523079302
$ws.Range("T5").Value = 0.183054523079302

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.103385
$ws.Range("H6").Value = 45.31015499999999
$ws.Range("I6").Value = 0.214656652056136
$ws.Range("J6").Value = 0.214656652056136
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.224107666666666
$ws.Range("N6").Value = 9.672322999999999
$ws.Range("O6").Value = 0.09060383010029813
$ws.Range("P6").Value = 0.09060383010029811
$ws.Range("Q6").Value = 48.69493937111832
$ws.Range("R6").Value = 438.2544543400649
$ws.Range("S6").Value = 0.01944871483279295
$ws.Range("T6").Value = 0.01944871483279295

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.103385
$ws.Range("H7").Value = 45.31015499999999
$ws.Range("I7").Value = 0.214656652056136
$ws.Range("J7").Value = 0.214656652056136
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.014730333333333
$ws.Range("N7").Value = 6.044191
$ws.Range("O7").Value = 0.05661792461415433
$ws.Range("P7").Value = 0.05661792461415433
$ws.Range("Q7").Value = 30.42924789551166
$ws.Range("R7").Value = 273.863231059605
$ws.Range("S7").Value = 0.01215341414404106
$ws.Range("T7").Value = 0.01215341414404106

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.62086333333334
$ws.Range("H8").Value = 142.86259
$ws.Range("I8").Value = 0.676810866647188
$ws.Range("J8").Value = 0.676810866647188
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 30.345835
$ws.Range("N8").Value = 91.03750500000001
$ws.Range("O8").Value = 0.8527782452855476
$ws.Range("P8").Value = 0.8527782452855475
$ws.Range("Q8").Value = 1445.094861270884
$ws.Range("R8").Value = 13005.85375143795
$ws.Range("S8").Value = 0.5771695832495798
$ws.Range("T8").Value = 0.5771695832495797

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.62086333333334
$ws.Range("H9").Value = 142.86259
$ws.Range("I9").Value = 0.676810866647188
$ws.Range("J9").Value = 0.676810866647188
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.224107666666666
$ws.Range("N9").Value = 9.672322999999999
$ws.Range("O9").Value = 0.09060383010029813
$ws.Range("P9").Value = 0.09060383010029811
$ws.Range("Q9").Value = 153.5347905662856
$ws.Range("R9").Value = 1381.81311509657
$ws.Range("S9").Value = 0.06132165677173736
$ws.Range("T9").Value = 0.06132165677173735

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 47.62086333333334
$ws.Range("H10").Value = 142.86259
$ws.Range("I10").Value = 0.676810866647188
$ws.Range("J10").Value = 0.676810866647188
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.014730333333333
$ws.Range("N10").Value = 6.044191
$ws.Range("O10").Value = 0.05661792461415433
$ws.Range("P10").Value = 0.05661792461415433
$ws.Range("Q10").Value = 95.94319785718778
$ws.Range("R10").Value = 863.48878071469
$ws.Range("S10").Value = 0.03831962662587095
$ws.Range("T10").Value = 0.03831962662587095
